{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: In the \"Copy source code DigitalMenu folder in \" paragraph,\n//   expand the text to mention both possible folder names/locations and\n//   add a trailing \"_GoBack\" bookmark (Word auto-drops one at the last\n//   edit location) right before \"C:\\inetpub\\wwwroot\".\n// Change 2: In the \"Permission errors e\" + bookmark + \"tc.\" paragraph,\n//   merge everything into a single run \"Permission errors etc.\" (the\n//   bookmark moved to change 1, so it disappears here).\n\nconst FLAT_OPC_HEADER =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>\\n' +\n  '<?mso-application progid=\"Word.Document\"?>\\n' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst FLAT_OPC_FOOTER =\n  '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>' +\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapParagraphOoxml(innerParagraphXml) {\n  return FLAT_OPC_HEADER + innerParagraphXml + FLAT_OPC_FOOTER;\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet copySourcePara = null;\nlet permissionPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (copySourcePara === null && t.indexOf(\"Copy source code\") !== -1 && t.indexOf(\"DigitalMenu\") !== -1) {\n    copySourcePara = paragraphs.items[i];\n  }\n  if (permissionPara === null && t.indexOf(\"Permission errors\") !== -1) {\n    permissionPara = paragraphs.items[i];\n  }\n}\n\nif (!copySourcePara) {\n  throw new Error('Could not find the \"Copy source code ... DigitalMenu folder in\" paragraph.');\n}\nif (!permissionPara) {\n  throw new Error('Could not find the \"Permission errors etc.\" paragraph.');\n}\n\n// --- Change 1 -----------------------------------------------------------\nconst copySourceXml = wrapParagraphOoxml(\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Copy source code </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>DigitalMenu</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> or </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>DigitalMenu</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> Inside of Published</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> folder in</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> to</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>C:\\\\inetpub\\\\wwwroot</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> folder.</w:t></w:r>' +\n  '</w:p>'\n);\ncopySourcePara.insertOoxml(copySourceXml, \"Replace\");\nawait context.sync();\n\n// --- Change 2 -------------------------------------------------------------\nconst permissionXml = wrapParagraphOoxml(\n  '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t>Permission errors etc.</w:t></w:r>' +\n  '</w:p>'\n);\npermissionPara.insertOoxml(permissionXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change 1: In the \"Copy source code DigitalMenu folder in \" paragraph,\n#   expand the text to mention both possible folder names/locations and\n#   add a trailing \"_GoBack\" bookmark (Word drops this marker at the last\n#   edit location) right before \"C:\\inetpub\\wwwroot\".\n# Change 2: In the \"Permission errors e\" + bookmark + \"tc.\" paragraph,\n#   merge everything into a single run \"Permission errors etc.\" (the\n#   bookmark moved to change 1, so it disappears here).\n\n$d = $word.ActiveDocument\n\n$flatOpcHeader = '<?xml version=\"1.0\" standalone=\"yes\"?>' + \"`n\" +\n  '<?mso-application progid=\"Word.Document\"?>' + \"`n\" +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>'\n$flatOpcFooter = '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfunction Wrap-ParagraphOoxml($innerParagraphXml) {\n    return $flatOpcHeader + $innerParagraphXml + $flatOpcFooter\n}\n\n# --- Locate the two target paragraphs by their text content --------------\n$copySourceIndex = -1\n$permissionIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($copySourceIndex -eq -1 -and $t -like \"*Copy source code*\" -and $t -like \"*DigitalMenu*\") {\n        $copySourceIndex = $i\n    }\n    if ($permissionIndex -eq -1 -and $t -like \"*Permission errors*\") {\n        $permissionIndex = $i\n    }\n}\n\nif ($copySourceIndex -eq -1) {\n    throw \"Could not find the 'Copy source code ... DigitalMenu folder in' paragraph.\"\n}\nif ($permissionIndex -eq -1) {\n    throw \"Could not find the 'Permission errors etc.' paragraph.\"\n}\n\n# --- Change 1 --------------------------------------------------------------\n$copySourceInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Copy source code </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>DigitalMenu</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> or </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>DigitalMenu</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> Inside of Published</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> folder in</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> to</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t>C:\\inetpub\\wwwroot</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> folder.</w:t></w:r>' +\n  '</w:p>'\n$copySourceRange = $d.Paragraphs.Item($copySourceIndex).Range\n[void]$copySourceRange.InsertXML((Wrap-ParagraphOoxml $copySourceInner))\n\n# --- Change 2 --------------------------------------------------------------\n# Re-fetch: paragraph count/order is unaffected by change 1 (same number of\n# paragraphs), so the previously located index is still valid.\n$permissionInner = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr></w:pPr>' +\n  '<w:r><w:t>Permission errors etc.</w:t></w:r>' +\n  '</w:p>'\n$permissionRange = $d.Paragraphs.Item($permissionIndex).Range\n[void]$permissionRange.InsertXML((Wrap-ParagraphOoxml $permissionInner))\n\nWrite-Output \"done\"\n"}
